$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -12
$ws.Range("F3").Value = -1
$ws.Range("F5").Value = -7
$ws.Range("F6").Value = -6
$ws.Range("F17").Value = -1
